# Add 7 new FOSS4Spectroscopy database entries (rows 127-133) and fix a couple
# of adjacent links, per commit "add to database; fix links".
#
# Cell values are written in the same scattered order the original author
# used (name/repo/description first, details filled in afterwards) so the
# resulting shared-string table matches the source edit exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values reused across several of the new rows -------------------------
$ws.Range("B127").Value = "NA"
$ws.Range("B128").Value = "NA"
$ws.Range("B129").Value = "NA"

$ws.Range("G131").Value = "Any"
$ws.Range("G132").Value = "Any"

$ws.Range("E127").Value = "Python"
$ws.Range("E129").Value = "Python"
$ws.Range("E130").Value = "Python"
$ws.Range("E131").Value = "Python"
$ws.Range("E132").Value = "Python"
$ws.Range("E133").Value = "Python"

$ws.Range("G127").Value = "NMR"
$ws.Range("G129").Value = "NMR"
$ws.Range("G133").Value = "NMR"

$ws.Range("G128").Value = "Data sharing (NMR nD)"

# --- Row 127: ANSURR --------------------------------------------------------
$ws.Range("A127").Value = "ANSURR"
$ws.Range("C127").Value = "https://github.com/nickjf/ANSURR2"
$ws.Range("F127").Value = "Validate accuracy of protein NMR structures"
$ws.Range("D127").Value = "https://www.nature.com/articles/s41467-020-20177-1"

# --- Row 128: nmr-parser ----------------------------------------------------
$ws.Range("A128").Value = "nmr-parser"
$ws.Range("C128").Value = "https://github.com/cheminfo/nmr-parser"
$ws.Range("E128").Value = "JavaScript"
$ws.Range("F128").Value = "Parse 1D & 2D NMR files to JSON"

# --- Row 129: nmrgnn ---------------------------------------------------------
$ws.Range("A129").Value = "nmrgnn"
$ws.Range("C129").Value = "https://github.com/ur-whitelab/nmrgnn"
$ws.Range("F129").Value = "Graph neural network prediction of NMR shifts"
$ws.Range("D129").Value = "https://pubs.rsc.org/en/content/articlehtml/2021/sc/d1sc01895g"

# --- Row 130: cwepr ----------------------------------------------------------
$ws.Range("A130").Value = "cwepr"
$ws.Range("B130").Value = "https://docs.cwepr.de/v0.2/"
$ws.Range("C130").Value = "https://github.com/tillbiskup/cwepr"
$ws.Range("F130").Value = "Processing continuous wave EPR data"
$ws.Range("G130").Value = "EPR"

# --- Row 131: ASpecD ---------------------------------------------------------
$ws.Range("A131").Value = "ASpecD"
$ws.Range("B131").Value = "https://www.aspecd.de/public/index"
$ws.Range("C131").Value = "https://github.com/tillbiskup/aspecd"
$ws.Range("F131").Value = "General handling of spectroscopic data"

# --- Row 132: trEPR -----------------------------------------------------------
$ws.Range("A132").Value = "trEPR"
$ws.Range("F132").Value = "Time-resolved EPR spectroscopy"
$ws.Range("C132").Value = "https://github.com/tillbiskup/trepr"
$ws.Range("B132").Value = "https://docs.trepr.de/v0.2/index.html"

# --- Row 133: NMR-EsPy --------------------------------------------------------
$ws.Range("A133").Value = "NMR-EsPy"
$ws.Range("D133").Value = "https://www.sciencedirect.com/science/article/pii/S1090780722000313"
$ws.Range("C133").Value = "https://github.com/foroozandehgroup/NMR-EsPy"
$ws.Range("B133").Value = "https://foroozandehgroup.github.io/NMR-EsPy/"
$ws.Range("F133").Value = "Estimation of NMR parameters"

# --- leave the view/selection the way the author left it --------------------
[void]$ws.Range("F134").Select()
